# Fill down the station name (column A) for every "false spring year" detail
# row on Sheet4 so each row carries its station label -- needed to build the
# new linear model of false-spring-event frequency per station.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")

$stationRows = @{
    "Anthony, KS"    = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14)
    "Hastings, NE"   = @(16, 17, 18, 19, 20, 21)
    "West Point, NE" = @(23, 24, 25, 26)
    "Yankton, SD"    = @(28, 29, 30, 31)
    "Aberdeen, SD"   = @(34)
}

foreach ($station in $stationRows.Keys) {
    foreach ($r in $stationRows[$station]) {
        $ws.Cells.Item($r, 1).Value = $station
    }
}

# Leave the selection where the edit finished, matching the author's last
# active cell on the (still) active Sheet4 tab.
$ws.Range("A34").Select()
